# Apply updated crypto price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.200.70"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "2.952.59"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("E6").Value = "  -1.41%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0838"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "3.423.16"
$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "2.950.53"
$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("E17").Value = "  +4.96%  "

$ws.Range("D18").Value = "51.276.89"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.63%  "

$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.30%  "

$ws.Range("E27").Value = "  +8.40%  "

$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.112"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.15%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0446"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.73%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  -1.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.47%  "

$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.271"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("D49").Value = "2.001.65"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("E50").Value = "  +1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0336"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.30%  "

